$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("choices")
$ws2.Columns.Item(1).ColumnWidth = 48.330729
$ws2.Columns.Item(2).ColumnWidth = 33.666667
$ws2.Columns.Item(3).ColumnWidth = 52.498698
